# Updates the cryptocurrency table on Sheet1 with refreshed price/volume data.
# A handful of "Price" values are plain decimal numbers (e.g. "576.17"); Excel's
# automatic type detection would otherwise store them as numbers and normalize
# their text (e.g. "1.00" -> 1, "0.940" -> 0.94). Prefixing such values with a
# leading apostrophe forces Excel to keep them as literal text, matching the
# original inline-string cell content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    if ($text -match "^[+-]?\d+(\.\d+)?$") {
        $ws.Range($cellRef).Value = "'" + $text
    } else {
        $ws.Range($cellRef).Value = $text
    }
}

Set-TextValue 'D2' '66.694.11'
Set-TextValue 'D3' '3.088.25'
Set-TextValue 'E3' '  -1.26%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '576.17'
Set-TextValue 'E5' '  -0.91%  '
Set-TextValue 'D6' '172.38'
Set-TextValue 'E6' '  -1.17%  '
Set-TextValue 'E7' '  +0.02%  '
Set-TextValue 'D8' '3.085.71'
Set-TextValue 'E8' '  -1.20%  '
Set-TextValue 'E9' '  -1.76%  '
Set-TextValue 'D10' '6.35'
Set-TextValue 'E10' '  -0.95%  '
Set-TextValue 'E11' '  -2.97%  '
Set-TextValue 'E12' '  -2.44%  '
Set-TextValue 'E13' '  -4.51%  '
Set-TextValue 'D14' '35.67'
Set-TextValue 'E14' '  -4.65%  '
Set-TextValue 'D15' '0.121'
Set-TextValue 'E15' '  -0.78%  '
Set-TextValue 'D16' '3.603.77'
Set-TextValue 'D17' '66.625.86'
Set-TextValue 'E17' '  -0.76%  '
Set-TextValue 'B18' 'Chainlink'
Set-TextValue 'C18' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D18' '16.91'
Set-TextValue 'E18' '  +2.63%  '
Set-TextValue 'B19' 'Polkadot'
Set-TextValue 'C19' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D19' '6.95'
Set-TextValue 'E19' '  -2.73%  '
Set-TextValue 'D20' '3.089.56'
Set-TextValue 'E20' '  -1.18%  '
Set-TextValue 'D21' '484.01'
Set-TextValue 'E21' '  -1.93%  '
Set-TextValue 'D22' '7.73'
Set-TextValue 'E22' '  -2.26%  '
Set-TextValue 'D23' '0.688'
Set-TextValue 'E23' '  -3.08%  '
Set-TextValue 'D24' '83.31'
Set-TextValue 'E24' '  -1.14%  '
Set-TextValue 'D25' '12.65'
Set-TextValue 'E25' '  -5.10%  '
Set-TextValue 'D26' '2.22'
Set-TextValue 'E26' '  -3.35%  '
Set-TextValue 'D27' '10.04'
Set-TextValue 'E27' '  -3.94%  '
Set-TextValue 'E28' '  -0.05%  '
Set-TextValue 'D29' '7.94'
Set-TextValue 'E29' '  +0.03%  '
Set-TextValue 'D30' '2.25'
Set-TextValue 'E30' '  -4.48%  '
Set-TextValue 'E31' '  -4.25%  '
Set-TextValue 'D32' '27.86'
Set-TextValue 'E32' '  -3.10%  '
Set-TextValue 'E33' '  -3.07%  '
Set-TextValue 'D34' '0.0₃0929'
Set-TextValue 'E34' '  -2.23%  '
Set-TextValue 'D35' '1.00'
Set-TextValue 'E35' '  +0.03%  '
Set-TextValue 'D36' '48.25'
Set-TextValue 'E36' '  +2.37%  '
Set-TextValue 'E37' '  -5.87%  '
Set-TextValue 'D38' '0.940'
Set-TextValue 'E38' '  -3.73%  '
Set-TextValue 'E39' '  -2.34%  '
Set-TextValue 'D40' '0.308'
Set-TextValue 'E40' '  -1.31%  '
Set-TextValue 'E41' '  -1.62%  '
Set-TextValue 'D42' '1.95'
Set-TextValue 'E42' '  -5.55%  '
Set-TextValue 'D43' '8.25'
Set-TextValue 'E43' '  -3.71%  '
Set-TextValue 'D44' '2.59'
Set-TextValue 'E44' '  -1.06%  '
Set-TextValue 'D45' '2.776.65'
Set-TextValue 'E45' '  -2.17%  '
Set-TextValue 'D46' '0.0345'
Set-TextValue 'E46' '  -2.44%  '
Set-TextValue 'D47' '366.88'
Set-TextValue 'E47' '  -4.82%  '
Set-TextValue 'D48' '134.28'
Set-TextValue 'E48' '  -1.14%  '
Set-TextValue 'D50' '24.37'
Set-TextValue 'E50' '  -2.67%  '
Set-TextValue 'E51' '  -2.83%  '
